$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.928.59"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "1.553.20"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.65"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.70"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D12").Value = "1.773.72"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "1.554.23"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "26.918.66"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.64"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.05"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.22"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.73"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.56"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "1.422.52"
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.519"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.67"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.986"
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.65"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "1.688.14"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.20"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0524"
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0956"
$ws.Range("E51").Value = "  +1.47%  "